# Add "The House of Morgan" to the Completed reading list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")
$ws.Activate()

# New row will be row 30 (right after the current last row, 29).
$newRow = 30

$ws.Range("A$newRow").Value = "The House of Morgan"
$ws.Range("B$newRow").Value = "Ron Chernow"

# Copy the date formatting/style from the row above, then overwrite the
# values, so the new date cells reuse the existing date style (same as the
# rest of the column) instead of creating a brand new number format.
$ws.Range("C29").Copy($ws.Range("C$newRow"))
$ws.Range("C$newRow").Value = 43881   # 2/20/2020 (Start Date)

$ws.Range("D29").Copy($ws.Range("D$newRow"))
$ws.Range("D$newRow").Value = 43888   # 2/27/2020 (Finish Date)

$ws.Range("E$newRow").Value = "finance;history;jp morgan;jack morgan;great depression;investment banking;history"
$ws.Range("F$newRow").Value = "Audio"
$ws.Range("G$newRow").Value = "34 Hours 37 Mins"

# Update the view/scroll state to reflect the newly added row, matching
# where Excel would land after typing the new entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("B31").Select() | Out-Null
